# [TEST] : modif fichier test
# Reproduces the commit: clears three stray "Date" (column B) cells that
# duplicated column A's timestamp, and appends four new log rows (14-17)
# to the "Reponses au formulaire 1" sheet, mirroring the formatting of
# existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Réponses au formulaire 1")

# ---------------------------------------------------------------------
# 1) Fully remove the now-unwanted column-B cells on rows 2, 3 and 9
#    (Clear() drops both value AND style, matching the cells disappearing
#    from the sheet entirely rather than merely becoming blank).
# ---------------------------------------------------------------------
$ws.Range("B2").Clear()
$ws.Range("B3").Clear()
$ws.Range("B9").Clear()

# ---------------------------------------------------------------------
# 2) Append four new rows of data (14-17). Each column's formatting is
#    copied from an existing row that already carries the right style,
#    then the value is overwritten - this reuses the workbook's existing
#    cellXfs entries instead of minting new ones.
# ---------------------------------------------------------------------
function Copy-CellValue {
    param($SrcAddr, $DstAddr, $Value)
    $ws.Range($SrcAddr).Copy() | Out-Null
    $ws.Range($DstAddr).PasteSpecial(-4122) | Out-Null
    $ws.Range($DstAddr).Value = $Value
}

# Row 14 (pattern mirrors row 6: arrival of "Louise, Joséphine")
Copy-CellValue "A6" "A14" 43356.868058912034
Copy-CellValue "B6" "B14" 43356
Copy-CellValue "C6" "C14" "Louise, Joséphine"
Copy-CellValue "D6" "D14" "Arrivée"
Copy-CellValue "E6" "E14" 0.32291666666424135

# Row 15 (pattern mirrors row 7: departure of "Joséphine" with "AR école")
Copy-CellValue "A7" "A15" 43356.868448819441
Copy-CellValue "B7" "B15" 43356
Copy-CellValue "C7" "C15" "Joséphine"
Copy-CellValue "D7" "D15" "Départ"
Copy-CellValue "E7" "E15" 0.72916666666424135
Copy-CellValue "G7" "G15" "AR école"

# Row 16 (pattern mirrors row 8: arrival of "Joséphine")
Copy-CellValue "A8" "A16" 43357.862597106483
Copy-CellValue "B8" "B16" 43357
Copy-CellValue "C8" "C16" "Joséphine"
Copy-CellValue "D8" "D16" "Arrivée"
Copy-CellValue "E8" "E16" 0.375

# Row 17 (pattern mirrors row 11: departure of "Joséphine")
Copy-CellValue "A11" "A17" 43357.862827986115
Copy-CellValue "B11" "B17" 43357
Copy-CellValue "C11" "C17" "Joséphine"
Copy-CellValue "D11" "D17" "Départ"
Copy-CellValue "E11" "E17" 0.70833333333575865

# Match the row height used throughout the rest of the sheet (15.75pt,
# explicit/custom) on the four freshly-added rows.
$ws.Rows("14:17").RowHeight = 15.75

# ---------------------------------------------------------------------
# 3) Move the active selection down to B18, just past the freshly added
#    rows (matches where the author's cursor ended up after typing).
# ---------------------------------------------------------------------
$ws.Range("B18").Select() | Out-Null
